# MAJ mapping suite review de NRISS 39b045f35aead0577830ffad59c0bbbf16e0bb22
#
# Applies the following changes to FRAllergyIntoleranceLMCDAFHIR.xlsx:
#  - Metadata sheet: clear the "Name" value (B4), replace the "Title" value (B5)
#    with the former "Name" text, and bump the "Date" value (B8).
#  - Mapping Table sheets: rename "entryRelationship." segments to "entryRelationship:"
#    in the four FRCDAAllergieOuHypersensibilite.entryRelationship.* target strings.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")

$nameValue = $metadata.Range("B4").Value2
$metadata.Range("B4").Value = ""
$metadata.Range("B5").Value = $nameValue
$metadata.Range("B8").Value = "2026-01-07T15:20:53+00:00"

$table0 = $wb.Worksheets.Item("Mapping Table 0")
$table0.Range("D10").Value = "FRCDAAllergieOuHypersensibilite.entryRelationship:frProbleme"
$table0.Range("D11").Value = "FRCDAAllergieOuHypersensibilite.entryRelationship:frStatutCliniqueAllergie"
$table0.Range("D12").Value = "FRCDAAllergieOuHypersensibilite.entryRelationship:frCertitude"
$table0.Range("D13").Value = "FRCDAAllergieOuHypersensibilite.entryRelationship:frCriticite"

$table1 = $wb.Worksheets.Item("Mapping Table 1")
$table1.Range("A8").Value = "FRCDAAllergieOuHypersensibilite.entryRelationship:frProbleme"
$table1.Range("A10").Value = "FRCDAAllergieOuHypersensibilite.entryRelationship:frStatutCliniqueAllergie"
$table1.Range("A11").Value = "FRCDAAllergieOuHypersensibilite.entryRelationship:frCertitude"
$table1.Range("A12").Value = "FRCDAAllergieOuHypersensibilite.entryRelationship:frCriticite"
